# "Correction de mini bug"
# The producer name ("La ferme de Robert") shown in B1 was wrong / stray
# for this bill, so it is cleared out (cell keeps its style, just no
# longer holds a value -> the shared string itself disappears too since
# it was only referenced once). A handful of the product-line quantities
# also get corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the producer-name cell (was "La ferme de Robert").
$ws.Range("B1").Value = ""

# Fix the quantities ordered for each line item.
$ws.Range("B9").Value  = 4   # Pain complet
$ws.Range("B14").Value = 2   # Tomates grappe
$ws.Range("B19").Value = 3   # Pomme de terre
$ws.Range("B24").Value = 1   # Salade
$ws.Range("B29").Value = 8   # Radis
